# activity_relation_data_preparation: add non-related relations for the
# remaining activity pairs.
#
# - Relabel the "gateway_1"/"gateway_2" headers as "activity_1"/"activity_2"
#   on every sheet that has them.
# - Rename "Doc Count" -> "Doc Count Normal" and "Doc Stats" -> "Doc Stats
#   Normal".
# - Append a "non_related" summary row to "Relation Type Count" and
#   "Comment Count".
# - Add two new sheets - "Doc Count Non-related" and "Doc Stats Non-related"
#   - mirroring the structure of the "Normal" sheets, for the non_related
#   activity pairs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Relation Type Count": header rename + new "non_related" row.
# ---------------------------------------------------------------------------
$wsRelType = $wb.Worksheets.Item("Relation Type Count")
$wsRelType.Range("C1").Value2 = "activity_1"
$wsRelType.Range("D1").Value2 = "activity_2"

$wsRelType.Range("A4").Copy()
$wsRelType.Range("A5").PasteSpecial(-4122)
$wsRelType.Range("A5").Value2 = "non_related"
$wsRelType.Range("B5:E5").Value2 = 2959

# ---------------------------------------------------------------------------
# 2. "Comment Count": header rename + new "non_related" row.
# ---------------------------------------------------------------------------
$wsComment = $wb.Worksheets.Item("Comment Count")
$wsComment.Range("D1").Value2 = "activity_1"
$wsComment.Range("E1").Value2 = "activity_2"

$wsComment.Range("A8").Copy()
$wsComment.Range("A9").PasteSpecial(-4122)
$wsComment.Range("A9").Value2 = "non_related"

$wsComment.Range("B8").Copy()
$wsComment.Range("B9").PasteSpecial(-4122)
$wsComment.Range("B9").ClearContents()

$wsComment.Range("C9:E9").Value2 = 2959

# ---------------------------------------------------------------------------
# 3. "Doc Count" -> "Doc Count Normal": header rename + sheet rename.
# ---------------------------------------------------------------------------
$wsDocCount = $wb.Worksheets.Item("Doc Count")
$wsDocCount.Range("B1").Value2 = "activity_1"
$wsDocCount.Range("C1").Value2 = "activity_2"
$wsDocCount.Name = "Doc Count Normal"

# ---------------------------------------------------------------------------
# 4. "Doc Stats" -> "Doc Stats Normal": header rename + sheet rename.
# ---------------------------------------------------------------------------
$wsDocStats = $wb.Worksheets.Item("Doc Stats")
$wsDocStats.Range("B1").Value2 = "activity_1"
$wsDocStats.Range("C1").Value2 = "activity_2"
$wsDocStats.Name = "Doc Stats Normal"

# ---------------------------------------------------------------------------
# 5. New sheet "Doc Count Non-related" (mirrors "Doc Count Normal").
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsDocCountNr = $wb.Worksheets.Add($null, $lastSheet)
$wsDocCountNr.Name = "Doc Count Non-related"

$wsDocCount.Range("A1:E1").Copy()
$wsDocCountNr.Range("A1:E1").PasteSpecial(-4122)
$wsDocCountNr.Range("A1").Value2 = "doc_name"
$wsDocCountNr.Range("B1").Value2 = "activity_1"
$wsDocCountNr.Range("C1").Value2 = "activity_2"
$wsDocCountNr.Range("D1").Value2 = "relation_type"
$wsDocCountNr.Range("E1").Value2 = "comment"

$wsDocCount.Range("A2:A45").Copy()
$wsDocCountNr.Range("A2:A45").PasteSpecial(-4122)

$wsDocCountNr.Range("A2").Value2 = "doc-1.1"
$wsDocCountNr.Range("B2:E2").Value2 = 33
$wsDocCountNr.Range("A3").Value2 = "doc-1.2"
$wsDocCountNr.Range("B3:E3").Value2 = 31
$wsDocCountNr.Range("A4").Value2 = "doc-1.3"
$wsDocCountNr.Range("B4:E4").Value2 = 42
$wsDocCountNr.Range("A5").Value2 = "doc-1.4"
$wsDocCountNr.Range("B5:E5").Value2 = 44
$wsDocCountNr.Range("A6").Value2 = "doc-10.1"
$wsDocCountNr.Range("B6:E6").Value2 = 2
$wsDocCountNr.Range("A7").Value2 = "doc-10.10"
$wsDocCountNr.Range("B7:E7").Value2 = 30
$wsDocCountNr.Range("A8").Value2 = "doc-10.11"
$wsDocCountNr.Range("B8:E8").Value2 = 18
$wsDocCountNr.Range("A9").Value2 = "doc-10.12"
$wsDocCountNr.Range("B9:E9").Value2 = 6
$wsDocCountNr.Range("A10").Value2 = "doc-10.13"
$wsDocCountNr.Range("B10:E10").Value2 = 1
$wsDocCountNr.Range("A11").Value2 = "doc-10.14"
$wsDocCountNr.Range("B11:E11").Value2 = 9
$wsDocCountNr.Range("A12").Value2 = "doc-10.2"
$wsDocCountNr.Range("B12:E12").Value2 = 72
$wsDocCountNr.Range("A13").Value2 = "doc-10.3"
$wsDocCountNr.Range("B13:E13").Value2 = 36
$wsDocCountNr.Range("A14").Value2 = "doc-10.4"
$wsDocCountNr.Range("B14:E14").Value2 = 35
$wsDocCountNr.Range("A15").Value2 = "doc-10.5"
$wsDocCountNr.Range("B15:E15").Value2 = 6
$wsDocCountNr.Range("A16").Value2 = "doc-10.6"
$wsDocCountNr.Range("B16:E16").Value2 = 3
$wsDocCountNr.Range("A17").Value2 = "doc-10.7"
$wsDocCountNr.Range("B17:E17").Value2 = 21
$wsDocCountNr.Range("A18").Value2 = "doc-10.8"
$wsDocCountNr.Range("B18:E18").Value2 = 20
$wsDocCountNr.Range("A19").Value2 = "doc-10.9"
$wsDocCountNr.Range("B19:E19").Value2 = 16
$wsDocCountNr.Range("A20").Value2 = "doc-2.1"
$wsDocCountNr.Range("B20:E20").Value2 = 744
$wsDocCountNr.Range("A21").Value2 = "doc-2.2"
$wsDocCountNr.Range("B21:E21").Value2 = 302
$wsDocCountNr.Range("A22").Value2 = "doc-3.1"
$wsDocCountNr.Range("B22:E22").Value2 = 45
$wsDocCountNr.Range("A23").Value2 = "doc-3.2"
$wsDocCountNr.Range("B23:E23").Value2 = 5
$wsDocCountNr.Range("A24").Value2 = "doc-3.3"
$wsDocCountNr.Range("B24:E24").Value2 = 12
$wsDocCountNr.Range("A25").Value2 = "doc-3.5"
$wsDocCountNr.Range("B25:E25").Value2 = 102
$wsDocCountNr.Range("A26").Value2 = "doc-3.6"
$wsDocCountNr.Range("B26:E26").Value2 = 21
$wsDocCountNr.Range("A27").Value2 = "doc-3.7"
$wsDocCountNr.Range("B27:E27").Value2 = 10
$wsDocCountNr.Range("A28").Value2 = "doc-3.8"
$wsDocCountNr.Range("B28:E28").Value2 = 33
$wsDocCountNr.Range("A29").Value2 = "doc-4.1"
$wsDocCountNr.Range("B29:E29").Value2 = 558
$wsDocCountNr.Range("A30").Value2 = "doc-5.1"
$wsDocCountNr.Range("B30:E30").Value2 = 5
$wsDocCountNr.Range("A31").Value2 = "doc-5.2"
$wsDocCountNr.Range("B31:E31").Value2 = 13
$wsDocCountNr.Range("A32").Value2 = "doc-5.3"
$wsDocCountNr.Range("B32:E32").Value2 = 86
$wsDocCountNr.Range("A33").Value2 = "doc-5.4"
$wsDocCountNr.Range("B33:E33").Value2 = 30
$wsDocCountNr.Range("A34").Value2 = "doc-6.1"
$wsDocCountNr.Range("B34:E34").Value2 = 349
$wsDocCountNr.Range("A35").Value2 = "doc-6.2"
$wsDocCountNr.Range("B35:E35").Value2 = 6
$wsDocCountNr.Range("A36").Value2 = "doc-6.3"
$wsDocCountNr.Range("B36:E36").Value2 = 27
$wsDocCountNr.Range("A37").Value2 = "doc-7.1"
$wsDocCountNr.Range("B37:E37").Value2 = 13
$wsDocCountNr.Range("A38").Value2 = "doc-8.1"
$wsDocCountNr.Range("B38:E38").Value2 = 6
$wsDocCountNr.Range("A39").Value2 = "doc-8.2"
$wsDocCountNr.Range("B39:E39").Value2 = 27
$wsDocCountNr.Range("A40").Value2 = "doc-8.3"
$wsDocCountNr.Range("B40:E40").Value2 = 8
$wsDocCountNr.Range("A41").Value2 = "doc-9.1"
$wsDocCountNr.Range("B41:E41").Value2 = 25
$wsDocCountNr.Range("A42").Value2 = "doc-9.2"
$wsDocCountNr.Range("B42:E42").Value2 = 15
$wsDocCountNr.Range("A43").Value2 = "doc-9.3"
$wsDocCountNr.Range("B43:E43").Value2 = 28
$wsDocCountNr.Range("A44").Value2 = "doc-9.4"
$wsDocCountNr.Range("B44:E44").Value2 = 36
$wsDocCountNr.Range("A45").Value2 = "doc-9.5"
$wsDocCountNr.Range("B45:E45").Value2 = 28

# ---------------------------------------------------------------------------
# 6. New sheet "Doc Stats Non-related" (mirrors "Doc Stats Normal").
# ---------------------------------------------------------------------------
$wsDocStatsNr = $wb.Worksheets.Add($null, $wsDocCountNr)
$wsDocStatsNr.Name = "Doc Stats Non-related"

$wsDocStats.Range("B1:E1").Copy()
$wsDocStatsNr.Range("B1:E1").PasteSpecial(-4122)
$wsDocStatsNr.Range("B1").Value2 = "activity_1"
$wsDocStatsNr.Range("C1").Value2 = "activity_2"
$wsDocStatsNr.Range("D1").Value2 = "relation_type"
$wsDocStatsNr.Range("E1").Value2 = "comment"

# Column A labels ("count", "mean", ..., "25%", ..., "max") are identical
# to "Doc Stats Normal", so copy format AND literal value straight from
# there - this sidesteps Excel's autoconversion of strings such as "25%"
# into a percentage number when assigned through .Value2.
$wsDocStats.Range("A2:A9").Copy()
$wsDocStatsNr.Range("A2:A9").PasteSpecial(-4122)
$wsDocStats.Range("A2:A9").Copy()
$wsDocStatsNr.Range("A2:A9").PasteSpecial(-4163)

$wsDocStatsNr.Range("B2:E2").Value2 = 44
$wsDocStatsNr.Range("B3:E3").Value2 = 67.25
$wsDocStatsNr.Range("B4:E4").Value2 = 146.5137702870312
$wsDocStatsNr.Range("B5:E5").Value2 = 1
$wsDocStatsNr.Range("B6:E6").Value2 = 9.75
$wsDocStatsNr.Range("B7:E7").Value2 = 26
$wsDocStatsNr.Range("B8:E8").Value2 = 36
$wsDocStatsNr.Range("B9:E9").Value2 = 744
